# Commit: "added data from rerunning samples with sequential leach"
# Adds new rows of sample/cement data into the "cements" sheet (sheet 2),
# interleaved among the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Step 1: insert blank rows at the correct final positions -------------
# Work from the bottom of the sheet upward so earlier row numbers stay valid
# while we insert. (Row numbers below refer to the sheet's ORIGINAL layout.)

# 2 new rows before original row 24 (FUWM 16 T 6)  -> final rows 36-37
$ws.Rows.Item(24).Resize(2).Insert()

# 2 new rows before original row 17 (FUWM 8.5 4)    -> final rows 27-28
$ws.Rows.Item(17).Resize(2).Insert()

# 2 new rows before original row 14 (FUWM 8.5 1)    -> final rows 22-23
$ws.Rows.Item(14).Resize(2).Insert()

# 2 new rows before original row 13 (FUWM 3.5 H 3)  -> final rows 19-20
$ws.Rows.Item(13).Resize(2).Insert()

# 4 new rows before original row 11 (FUWM 3.5 H 1)  -> final rows 13-16
$ws.Rows.Item(11).Resize(4).Insert()

# 1 new row before original row 7 (FUWM 3.5 1)      -> final row 8
$ws.Rows.Item(7).Insert()

# 1 new row before original row 5 (FUWM 1 5)        -> final row 5
$ws.Rows.Item(5).Insert()

# --- Step 2: populate the new cells ----------------------------------------
# Column A (sample names) is written first, in the same order the samples
# were originally entered, so brand-new shared strings land in that order;
# column B (cement types) is filled in afterwards.

# "dil" (sequential-leach dilution) re-runs
$ws.Cells.Item(5, 1).Value  = "FUWM 1 3 dil"
$ws.Cells.Item(19, 1).Value = "FUWM 3.5 H 1 dil"
$ws.Cells.Item(20, 1).Value = "FUWM 3.5 H 2 dil"
$ws.Cells.Item(27, 1).Value = "FUWM 8.5 2 dil"
$ws.Cells.Item(28, 1).Value = "FUWM 8.5 3 dil"
$ws.Cells.Item(36, 1).Value = "FUWM 16 T 4 dil"
$ws.Cells.Item(37, 1).Value = "FUWM 16 T 5 dil"

# "T" transect samples
$ws.Cells.Item(13, 1).Value = "FUWM 3.5 T 1"
$ws.Cells.Item(14, 1).Value = "FUWM 3.5 T 2"
$ws.Cells.Item(15, 1).Value = "FUWM 3.5 T 3"
$ws.Cells.Item(16, 1).Value = "FUWM 3.5 T 5"

# new "3 A 1" sample
$ws.Cells.Item(8, 1).Value = "FUWM 3 A 1"
$ws.Cells.Item(8, 2).Value = "Micrite Above"

# new "4" samples
$ws.Cells.Item(22, 1).Value = "FUWM 4 1.5m B"
$ws.Cells.Item(23, 1).Value = "FUWM 4 3m B"
$ws.Cells.Item(22, 2).Value = "Micrite Below"

# Remaining column B (cement-type) values for the new rows, reusing the
# existing cement-type strings already present in the sheet.
$ws.Cells.Item(5, 2).Value  = "Micrite"
$ws.Cells.Item(13, 2).Value = "Micrite (SSF)"
$ws.Cells.Item(14, 2).Value = "Laminar Calcrete"
$ws.Cells.Item(15, 2).Value = "Laminar Microbial"
$ws.Cells.Item(16, 2).Value = "Laminar Calcrete"
$ws.Cells.Item(19, 2).Value = "Laminar Microbial"
$ws.Cells.Item(20, 2).Value = "Laminar Microbial"
$ws.Cells.Item(23, 2).Value = "Micrite Below"
$ws.Cells.Item(27, 2).Value = "Manganese"
$ws.Cells.Item(28, 2).Value = "Manganese"
$ws.Cells.Item(36, 2).Value = "Microbial"
$ws.Cells.Item(37, 2).Value = "Laminar Microbial"

# --- Step 3: cosmetic sheet/view updates to match the author's session -----
$ws.Columns.Item(1).ColumnWidth = 19.6
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D23").Select()
